# Update profit files after running on 2025-10-13
# Append the new day's row (date, profit) to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the date cell as Text first so Excel stores the literal string
# "10/13/2025" instead of auto-converting it to a date serial number
# (matching how the rest of the Date column is stored). ClearFormats()
# afterwards drops the temporary Text number-format so the cell is left
# with no explicit style, just like its neighbours.
$dateCell = $ws.Range("A57")
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/13/2025"
$dateCell.ClearFormats()

$ws.Range("B57").Value = 11910.66
